$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking strings as literal text
# (the source file has them as inlineStr, not numbers -- e.g. some keep
# trailing zeros like "1.000", others use a "." as a thousands separator
# like "29.533.33"). A plain .Value assignment of a numeric-looking string
# gets auto-converted by Excel into a real Number (dropping formatting, e.g.
# "1.000" -> 1, "120.00" -> 120, "0.000009973" -> 9.973E-06). Force those
# cells to Text first so the literal string is preserved exactly.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.533.33'
$ws.Range("E2").Value = '  +2.14%  '
$ws.Range("D3").Value = '1.990.69'
$ws.Range("E3").Value = '  +5.86%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '325.86'
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").Value = '0.4671'
$ws.Range("E7").Value = '  +1.39%  '
$ws.Range("D8").Value = '0.3945'
$ws.Range("E8").Value = '  +1.87%  '
$ws.Range("D9").Value = '46.41'
$ws.Range("E9").Value = '  -0.47%  '
$ws.Range("D10").Value = '0.07927'
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("D11").Value = '0.9999'
$ws.Range("E11").Value = '  +1.39%  '
$ws.Range("D12").Value = '22.87'
$ws.Range("E12").Value = '  +5.12%  '
$ws.Range("D13").Value = '2.025.89'
$ws.Range("E13").Value = '  +10.83%  '
$ws.Range("D14").Value = '7.261'
$ws.Range("E14").Value = '  +3.87%  '
$ws.Range("D15").Value = '5.868'
$ws.Range("E15").Value = '  +3.87%  '
$ws.Range("D16").Value = '0.07120'
$ws.Range("E16").Value = '  +2.30%  '
$ws.Range("D17").Value = '88.67'
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").Value = '0.000009973'
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").Value = '17.38'
$ws.Range("E20").Value = '  +2.55%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("D22").Value = '29.565.45'
$ws.Range("E22").Value = '  +2.27%  '
$ws.Range("D23").Value = '5.528'
$ws.Range("E23").Value = '  +5.59%  '
$ws.Range("D24").Value = '11.25'
$ws.Range("E24").Value = '  +2.72%  '
$ws.Range("D25").Value = '2.099'
$ws.Range("E25").Value = '  +0.60%  '
$ws.Range("D26").Value = '157.83'
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").Value = '19.61'
$ws.Range("E27").Value = '  +1.39%  '
$ws.Range("D28").Value = '5.993'
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = '120.00'
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").Value = '1.951'
$ws.Range("E30").Value = '  +1.45%  '
$ws.Range("D31").Value = '0.09434'
$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("D32").Value = '0.9084'
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("D33").Value = '1.349'
$ws.Range("E33").Value = '  +2.44%  '
$ws.Range("D34").Value = '5.250'
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("D35").Value = '3.179'
$ws.Range("E35").Value = '  -2.32%  '
$ws.Range("D36").Value = '0.000003495'
$ws.Range("E36").Value = '  +112.70%  '
$ws.Range("D37").Value = '0.05829'
$ws.Range("E37").Value = '  +1.48%  '
$ws.Range("E38").Value = '  -0.64%  '
$ws.Range("D39").Value = '0.02116'
$ws.Range("E39").Value = '  +2.14%  '
$ws.Range("D40").Value = '7.862'
$ws.Range("E40").Value = '  +3.12%  '
$ws.Range("D41").Value = '0.5746'
$ws.Range("E41").Value = '  +1.80%  '
$ws.Range("D42").Value = '0.1825'
$ws.Range("E42").Value = '  +3.44%  '
$ws.Range("D43").Value = '9.819'
$ws.Range("E43").Value = '  +1.39%  '
$ws.Range("D44").Value = '12.06'
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").Value = '2.715'
$ws.Range("E45").Value = '  +7.34%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.5367'
$ws.Range("E46").Value = '  +0.44%  '
$ws.Range("E47").Value = '  -5.06%  '
$ws.Range("D48").Value = '1.867'
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("D49").Value = '0.06944'
$ws.Range("E49").Value = '  -1.38%  '
$ws.Range("D50").Value = '114.18'
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("D51").Value = '0.3090'
$ws.Range("E51").Value = '  +7.96%  '
